# spain_laliga2_2023-2024.xlsx update
# 1) Three pairs of rows had their data (columns F:V) swapped between them.
# 2) Eight brand-new match rows were appended after the previous last row (68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Three pairs of rows have their match data (columns F:V) swapped between
# them (columns A:E - Indice/pais/torneio/temporada/data_partida - are
# identical between each swapped pair, so they are left untouched). Values
# are written explicitly (rather than read-swapped) for reliability.
# ---------------------------------------------------------------------------

# Rows 39 / 40 swap (Elche-Racing Santander  <->  Villarreal B-FC Cartagena SAD)
$ws.Range("F39").Value = "Villarreal B"
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = "FC Cartagena SAD"
$ws.Range("I39").Value = 2
$ws.Range("J39").Value = 2.37
$ws.Range("K39").Value = "26/08/2023 20:42"
$ws.Range("L39").Value = 2.11
$ws.Range("M39").Value = "02/09/2023 20:56"
$ws.Range("N39").Value = 3.33
$ws.Range("O39").Value = "26/08/2023 20:42"
$ws.Range("P39").Value = 3.47
$ws.Range("Q39").Value = "02/09/2023 20:52"
$ws.Range("R39").Value = 3.15
$ws.Range("S39").Value = "26/08/2023 20:42"
$ws.Range("T39").Value = 3.75
$ws.Range("U39").Value = "02/09/2023 20:59"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/spain/laliga2/villarreal-fc-cartagena-sad/2FXXUSq4/"

$ws.Range("F40").Value = "Elche"
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = "Racing Santander"
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 1.65
$ws.Range("K40").Value = "28/08/2023 23:42"
$ws.Range("L40").Value = 1.92
$ws.Range("M40").Value = "02/09/2023 20:57"
$ws.Range("N40").Value = 3.82
$ws.Range("O40").Value = "28/08/2023 23:42"
$ws.Range("P40").Value = 3.51
$ws.Range("Q40").Value = "02/09/2023 20:57"
$ws.Range("R40").Value = 5.66
$ws.Range("S40").Value = "28/08/2023 23:42"
$ws.Range("T40").Value = 4.46
$ws.Range("U40").Value = "02/09/2023 20:57"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/spain/laliga2/elche-racing-santander/CWkCDRET/"

# Rows 51 / 52 swap (Racing Santander-Amorebieta  <->  Ferrol-Villarreal B)
$ws.Range("F51").Value = "Ferrol"
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = "Villarreal B"
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = 2.22
$ws.Range("K51").Value = "04/09/2023 11:12"
$ws.Range("L51").Value = 1.88
$ws.Range("M51").Value = "10/09/2023 16:08"
$ws.Range("N51").Value = 3.13
$ws.Range("O51").Value = "04/09/2023 11:12"
$ws.Range("P51").Value = 3.53
$ws.Range("Q51").Value = "10/09/2023 16:08"
$ws.Range("R51").Value = 3.71
$ws.Range("S51").Value = "04/09/2023 11:12"
$ws.Range("T51").Value = 4.63
$ws.Range("U51").Value = "10/09/2023 16:08"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/spain/laliga2/ferrol-villarreal/lzCGM4ip/"

$ws.Range("F52").Value = "Racing Santander"
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = "Amorebieta"
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1.85
$ws.Range("K52").Value = "04/09/2023 11:12"
$ws.Range("L52").Value = 1.8
$ws.Range("M52").Value = "10/09/2023 16:14"
$ws.Range("N52").Value = 3.44
$ws.Range("O52").Value = "04/09/2023 11:12"
$ws.Range("P52").Value = 3.47
$ws.Range("Q52").Value = "10/09/2023 16:14"
$ws.Range("R52").Value = 4.79
$ws.Range("S52").Value = "04/09/2023 11:12"
$ws.Range("T52").Value = 5.36
$ws.Range("U52").Value = "10/09/2023 16:14"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/spain/laliga2/racing-santander-amorebieta/AucBrnE3/"

# Rows 53 / 54 swap (Mirandes-Andorra  <->  FC Cartagena SAD-Zaragoza)
$ws.Range("F53").Value = "FC Cartagena SAD"
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = "Zaragoza"
$ws.Range("I53").Value = 3
$ws.Range("J53").Value = 2.76
$ws.Range("K53").Value = "03/09/2023 17:43"
$ws.Range("L53").Value = 3.1
$ws.Range("M53").Value = "10/09/2023 18:23"
$ws.Range("N53").Value = 3.17
$ws.Range("O53").Value = "03/09/2023 17:43"
$ws.Range("P53").Value = 2.96
$ws.Range("Q53").Value = "10/09/2023 17:10"
$ws.Range("R53").Value = 2.77
$ws.Range("S53").Value = "03/09/2023 17:43"
$ws.Range("T53").Value = 2.7
$ws.Range("U53").Value = "10/09/2023 18:22"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-zaragoza/WfPlRpqT/"

$ws.Range("F54").Value = "Mirandes"
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = "Andorra"
$ws.Range("I54").Value = 3
$ws.Range("J54").Value = 2.85
$ws.Range("K54").Value = "03/09/2023 20:12"
$ws.Range("L54").Value = 2.88
$ws.Range("M54").Value = "10/09/2023 18:23"
$ws.Range("N54").Value = 3.16
$ws.Range("O54").Value = "03/09/2023 20:12"
$ws.Range("P54").Value = 3.12
$ws.Range("Q54").Value = "10/09/2023 17:08"
$ws.Range("R54").Value = 2.75
$ws.Range("S54").Value = "03/09/2023 20:12"
$ws.Range("T54").Value = 2.79
$ws.Range("U54").Value = "10/09/2023 18:27"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/spain/laliga2/mirandes-fc-andorra/xUmOAStc/"

# Rows 59 / 60 swap (Alcorcon-Levante  <->  Valladolid-FC Cartagena SAD)
$ws.Range("F59").Value = "Valladolid"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = "FC Cartagena SAD"
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1.93
$ws.Range("K59").Value = "13/09/2023 08:25"
$ws.Range("L59").Value = 1.78
$ws.Range("M59").Value = "16/09/2023 18:24"
$ws.Range("N59").Value = 3.57
$ws.Range("O59").Value = "13/09/2023 08:25"
$ws.Range("P59").Value = 3.71
$ws.Range("Q59").Value = "16/09/2023 18:24"
$ws.Range("R59").Value = 4.16
$ws.Range("S59").Value = "13/09/2023 08:25"
$ws.Range("T59").Value = 5.02
$ws.Range("U59").Value = "16/09/2023 18:27"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/spain/laliga2/valladolid-fc-cartagena-sad/lhXkfqDF/"

$ws.Range("F60").Value = "Alcorcon"
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = "Levante"
$ws.Range("I60").Value = 2
$ws.Range("J60").Value = 3.14
$ws.Range("K60").Value = "11/09/2023 20:13"
$ws.Range("L60").Value = 2.46
$ws.Range("M60").Value = "16/09/2023 18:29"
$ws.Range("N60").Value = 3.07
$ws.Range("O60").Value = "11/09/2023 20:13"
$ws.Range("P60").Value = 2.96
$ws.Range("Q60").Value = "16/09/2023 18:28"
$ws.Range("R60").Value = 2.53
$ws.Range("S60").Value = "11/09/2023 20:13"
$ws.Range("T60").Value = 3.53
$ws.Range("U60").Value = "16/09/2023 18:29"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/spain/laliga2/alcorcon-levante/tUzfzOLq/"

# ---------------------------------------------------------------------------
# Append 8 new rows (69-76) after the previous last row (68), cloning the
# formatting of row 68 (bold/bordered Indice column + date-formatted E
# column) down across the new range before filling in the values.
# ---------------------------------------------------------------------------

$ws.Range("A68:V68").Copy()
$ws.Range("A69:V76").PasteSpecial(-4122)

$newRows = @(
    @{ Row=69; A=68; E=45192.58333333334; F="Andorra"; G=0; H="Gijon"; I=0;
       J=1.98; K="18/09/2023 20:13"; L=2.31; M="23/09/2023 13:54";
       N=3.29; O="18/09/2023 20:13"; P=3.04; Q="23/09/2023 13:47";
       R=4.32; S="18/09/2023 20:13"; T=3.75; U="23/09/2023 13:54";
       V="https://www.betexplorer.com/football/spain/laliga2/fc-andorra-gijon/pfiW86B9/" },
    @{ Row=70; A=69; E=45192.67708333334; F="Levante"; G=2; H="Eldense"; I=0;
       J=1.76; K="18/09/2023 11:42"; L=1.75; M="23/09/2023 16:08";
       N=3.56; O="18/09/2023 11:42"; P=3.61; Q="23/09/2023 16:08";
       R=5.16; S="18/09/2023 11:42"; T=5.47; U="23/09/2023 16:08";
       V="https://www.betexplorer.com/football/spain/laliga2/levante-eldense/f5vdjzCs/" },
    @{ Row=71; A=70; E=45192.77083333334; F="Burgos CF"; G=4; H="Elche"; I=0;
       J=2.82; K="18/09/2023 11:42"; L=3.02; M="23/09/2023 18:29";
       N=3.13; O="18/09/2023 11:42"; P=3.05; Q="23/09/2023 17:50";
       R=2.75; S="18/09/2023 11:42"; T=2.71; U="23/09/2023 18:29";
       V="https://www.betexplorer.com/football/spain/laliga2/burgos-cf-elche/0pWTn15k/" },
    @{ Row=72; A=71; E=45192.875; F="Alcorcon"; G=0; H="Huesca"; I=2;
       J=2.4; K="18/09/2023 11:42"; L=2.31; M="23/09/2023 20:55";
       N=2.95; O="18/09/2023 11:42"; P=2.91; Q="23/09/2023 20:56";
       R=3.54; S="18/09/2023 11:42"; T=3.98; U="23/09/2023 20:56";
       V="https://www.betexplorer.com/football/spain/laliga2/alcorcon-huesca/2aDQmsjq/" },
    @{ Row=73; A=72; E=45192.875; F="Racing Santander"; G=2; H="Albacete"; I=1;
       J=2.56; K="16/09/2023 20:12"; L=2.66; M="23/09/2023 20:57";
       N=3.12; O="16/09/2023 20:12"; P=3; Q="23/09/2023 20:57";
       R=3.12; S="16/09/2023 20:12"; T=3.15; U="23/09/2023 20:57";
       V="https://www.betexplorer.com/football/spain/laliga2/racing-santander-albacete/l6CO8rS8/" },
    @{ Row=74; A=73; E=45193.58333333334; F="Villarreal B"; G=3; H="Amorebieta"; I=1;
       J=1.8; K="18/09/2023 11:42"; L=2.33; M="24/09/2023 13:56";
       N=3.71; O="18/09/2023 11:42"; P=3.35; Q="24/09/2023 13:56";
       R=4.62; S="18/09/2023 11:42"; T=3.34; U="24/09/2023 13:56";
       V="https://www.betexplorer.com/football/spain/laliga2/villarreal-amorebieta/rTBS72sF/" },
    @{ Row=75; A=74; E=45193.67708333334; F="R. Oviedo"; G=0; H="Valladolid"; I=1;
       J=2.74; K="18/09/2023 20:13"; L=2.51; M="24/09/2023 16:03";
       N=2.92; O="18/09/2023 20:13"; P=2.92; Q="24/09/2023 16:03";
       R=3.02; S="18/09/2023 20:13"; T=3.48; U="24/09/2023 16:03";
       V="https://www.betexplorer.com/football/spain/laliga2/r-oviedo-valladolid/AFDK9OC2/" },
    @{ Row=76; A=75; E=45193.77083333334; F="Mirandes"; G=1; H="Leganes"; I=3;
       J=2.47; K="17/09/2023 17:42"; L=2.29; M="24/09/2023 18:24";
       N=3.07; O="17/09/2023 17:42"; P=2.99; Q="24/09/2023 18:29";
       R=3.25; S="17/09/2023 17:42"; T=3.89; U="24/09/2023 18:20";
       V="https://www.betexplorer.com/football/spain/laliga2/mirandes-leganes/GKS4lEtf/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "spain"
    $ws.Range("C$row").Value = "laliga2"
    $ws.Range("D$row").Value = "2023-2024"
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
}

Write-Output "edit complete"
